$d = $word.ActiveDocument

# ± character (U+00B1) used in the margin-of-error bullet
$PM = [char]0xB1
$DOLLAR = "$"
# Hybrid bold + color (#2C3E50) highlight color as a Word BGR decimal value
$COLOR = 5258796

function Highlight-InParagraph($para, $text) {
    $r = $para.Range
    $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Font.Bold = $true
    $r.Font.Color = $COLOR
}

function Find-ParagraphContaining($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    return $null
}

# 1) "• Discovered systematic race coding errors ... from 23% to 64%"
$p1 = Find-ParagraphContaining("Discovered systematic race coding")
Highlight-InParagraph $p1 "23%"
Highlight-InParagraph $p1 "64%"

# 2) "• Utilized advanced sampling methods ... margin of error from ±4.2% to ±2.1%, ... from 71% to 87% ..."
$p2 = Find-ParagraphContaining("margin of error")
$pm1 = $PM + "4.2%"
$pm2 = $PM + "2.1%"
Highlight-InParagraph $p2 $pm1
Highlight-InParagraph $p2 $pm2
Highlight-InParagraph $p2 "71%"
Highlight-InParagraph $p2 "87%"

# 3) "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving ... $4.7M ..."
$p3 = Find-ParagraphContaining("Trigonometric algorithm for boundary")
$amt47 = $DOLLAR + "4.7M"
Highlight-InParagraph $p3 "73.5%"
Highlight-InParagraph $p3 $amt47

# 4) "• Built real-time FEC analysis systems ... valued over $2 trillion"
$p4 = Find-ParagraphContaining("Built real-time FEC")
$amt2 = $DOLLAR + "2"
Highlight-InParagraph $p4 $amt2

# 5) "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
$p5 = Find-ParagraphContaining("Algorithmic innovation")
Highlight-InParagraph $p5 "73.5%"

# 6) "• $4.7M savings enabled nonprofit access"
$p6 = Find-ParagraphContaining("savings enabled nonprofit")
$amt47b = $DOLLAR + "4.7M"
Highlight-InParagraph $p6 $amt47b

# 7) "• 178% accuracy improvement in racial classification algorithms"
$p7 = Find-ParagraphContaining("accuracy improvement in racial")
Highlight-InParagraph $p7 "178%"

Write-Output "done"
